$d = $word.ActiveDocument

# Step 1: merge the five separate runs ("в ____", " ", "часов _______", " ", "минут ")
# into a single run for BOTH occurrences in the document. The replacement text is
# identical to the concatenation of the original runs, so no wording changes -
# only the run layout is consolidated.
$d.Content.Find.Execute("в ____ часов _______ минут ", $true, $false, $false, $false, $false, $true, 1, $false, "в ____ часов _______ минут ", 2)

# Step 2: shorten the two distinct underscore blanks that follow "года" in each
# of the two occurrences (each has its own amount of leading whitespace, so they
# are addressed individually).
$d.Content.Find.Execute("года                      ______________", $true, $false, $false, $false, $false, $true, 1, $false, "года                      ___________", 2)
$d.Content.Find.Execute("года                       _____________", $true, $false, $false, $false, $false, $true, 1, $false, "года                       __________", 2)
